$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Carrera" column: rename Mecatronica -> Sistemas for both data rows
$ws.Range("B2").Value = "Sistemas"
$ws.Range("B3").Value = "Sistemas"

# Update the active selection to match the edited workbook's cursor position
$ws.Range("D6").Select()
